$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Set the "Private" column (D) to TRUE for every data row (2-6)
$ws.Range("D2:D6").Value = $true

# Column D (Private) now gets the same TRUE/FALSE list validation that
# column F (View) already carries, so rebuild the validations: the plain
# "any value" one on F1, the TRUE/FALSE list covering the rest of column F,
# and that same list validation newly applied to D2:D6.
$ws.Cells.Validation.Delete()
$ws.Range("F1").Validation.Add(0, 1, 1)
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("D2:D6").Validation.Add(3, 1, 1, """TRUE,FALSE""")

# Match the resulting selection state recorded in the saved file
$ws.Range("D2:D6").Select()
